$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level metadata -------------------------------------------------
# absolute path shown in the (hidden) x15ac:absPath element
$wb.Path = 'C:\Users\Roman\Dropbox\Projekte\Leistungsmessung\Breakout_Board_A1\Doku\'
# window size
$excel.Width = 28800
$excel.Height = 12435

# --- New header block (columns G:N, rows 1-2) --------------------------------
# Row 1: merged channel headers, bold + centered
$ws.Range("G1:H1").Merge()
$ws.Range("G1").Value = "Messung Kanal 1"

$ws.Range("I1:J1").Merge()
$ws.Range("I1").Value = "Messung Kanal 2"

$ws.Range("K1:L1").Merge()
$ws.Range("K1").Value = "Messung Kanal 3"

$ws.Range("M1:N1").Merge()
$ws.Range("M1").Value = "Messung Kanal 4"

$ws.Range("G1:N1").Font.Bold = $true
$ws.Range("G1:N1").HorizontalAlignment = -4108  ## xlCenter

# Row 2: sub-headers repeated for each channel
$ws.Range("G2").Value = "Strom sekundär in mA RMS"
$ws.Range("H2").Value = "Spannung sekundär in mA RMS"
$ws.Range("I2").Value = "Strom sekundär in mA RMS"
$ws.Range("J2").Value = "Spannung sekundär in mA RMS"
$ws.Range("K2").Value = "Strom sekundär in mA RMS"
$ws.Range("L2").Value = "Spannung sekundär in mA RMS"
$ws.Range("M2").Value = "Strom sekundär in mA RMS"
$ws.Range("N2").Value = "Spannung sekundär in mA RMS"

# --- Calibration measurement rows 3-6 ----------------------------------------
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 50.3

$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 55.6

$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 76.6

$ws.Range("G6").Value = 80
$ws.Range("H6").Value = 419
$ws.Range("I6").Value = 80
$ws.Range("J6").Value = 422
$ws.Range("K6").Value = 80
$ws.Range("L6").Value = 419.2
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 420.2

# --- Row 7: Mittelwert / Abweichung labels + formulas ------------------------
$ws.Range("F7").Value = "Mittelwert"
$ws.Range("G7").Value = "Abweichung in %"
$ws.Range("H7").Formula = '=((H6-$F$8)/$F$8)*100'
$ws.Range("I7").Value = "Abweichung"
$ws.Range("J7").Formula = '=((J6-$F$8)/$F$8)*100'
$ws.Range("K7").Value = "Abweichung"
$ws.Range("L7").Formula = '=((L6-$F$8)/$F$8)*100'
$ws.Range("M7").Value = "Abweichung"
$ws.Range("N7").Formula = '=((N6-$F$8)/$F$8)*100'

$ws.Range("H7,J7,L7,N7").NumberFormat = "0.0000"

# --- Row 8: average formula ---------------------------------------------------
$ws.Range("F8").Formula = "=AVERAGE(H6,J6,L6,N6)"

# --- Updated input value in row 20 (ripples through formulas in G20:L20) -----
$ws.Range("F20").Value = 35

# --- Column widths for the new columns ---------------------------------------
$ws.Range("G1").ColumnWidth = 25
$ws.Range("H1").ColumnWidth = 28.5703125
$ws.Range("I1").ColumnWidth = 25
$ws.Range("J1").ColumnWidth = 28.5703125
$ws.Range("K1").ColumnWidth = 25
$ws.Range("L1").ColumnWidth = 28.5703125
$ws.Range("M1").ColumnWidth = 25
$ws.Range("N1").ColumnWidth = 28.5703125

# --- Sheet view ----------------------------------------------------------------
$ws.Range("H9").Select()
$excel.ActiveWindow.ScrollColumn = 5

# --- Page setup ------------------------------------------------------------
$ws.PageSetup.PaperSize = 9      ## xlPaperA4
$ws.PageSetup.Orientation = 1    ## xlPortrait
